$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.204.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.59%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.515.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.14%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.31%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.16%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.51%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.07%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "

# Row 13
$ws.Range("E13").Value = "  -0.37%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.902.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.527.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.05%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.851"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.375.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.49%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0953"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.64%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "248.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.23%  "

# Row 25
$ws.Range("E25").Value = "  -3.63%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.42%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.58%  "

# Row 28
$ws.Range("E28").Value = "  +2.33%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.58%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.00%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.90%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.34%  "

# Row 34
$ws.Range("E34").Value = "  +0.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.98%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.66%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0776"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.91%  "

# Row 38
$ws.Range("E38").Value = "  -3.82%  "

# Row 39
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.14%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.15%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.03%  "

# Row 44
$ws.Range("E44").Value = "  -4.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0297"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.16%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.014.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.41%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "83.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.57%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.758.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.14%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "100.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
